$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the row above (row 6) into the new row 7 so both
# cells pick up the existing label/value column styling (fill + border),
# then set the new row's text.
$ws.Range("A6:B6").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)
$ws.Range("A7").Value = "email"

# Add the new value cell as a mailto: hyperlink (creates the built-in
# "Hyperlink" cell style + font, same as Excel does automatically). The
# display text is set on the cell itself first so Hyperlinks.Add doesn't
# need (and doesn't store) a separate TextToDisplay.
$ws.Range("B7").Value = "email1@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:email1@gmail.com")

# Match the row height used by the other rows.
$ws.Range("A7:B7").RowHeight = 39

# Leave the selection where the author left it.
[void]$ws.Range("D6").Select()
